$d = $word.ActiveDocument
$t = $d.Tables(1)

# NOTE: Find.Execute with Replace=2 (wdReplaceAll) searches/replaces across the
# *whole document*, even when called on a narrowly scoped Range - so we always
# use Replace=1 (wdReplaceOne) on a tightly scoped Range to touch only the
# single intended occurrence.

# ---- Row 3 (IMESD): 11-2-2017 -> 11-1-2018 ----
$t.Cell(3,1).Range.Find.Execute("11-2-2017", $true, $false, $false, $false, $false, $true, 1, $false, "11-1-2018", 1) | Out-Null
$t.Cell(3,2).Range.Paragraphs(1).Range.Find.Execute("Team: Brad Lenhardt, Gerald Tindal, & Sevrina Tindal", $true, $false, $false, $false, $false, $true, 1, $false, "Team: Brad Lenhardt, Brock Rowley", 1) | Out-Null

# ---- Row 4 (HDESD): 11-7-2017 -> 11-6-2018 ----
$t.Cell(4,1).Range.Find.Execute("11-7-2017", $true, $false, $false, $false, $false, $true, 1, $false, "11-6-2018", 1) | Out-Null
$t.Cell(4,2).Range.Paragraphs(1).Range.Find.Execute("Team: Brad Lenhardt, Gerald Tindal, & Sevrina Tindal", $true, $false, $false, $false, $false, $true, 1, $false, "Team: Brad Lenhardt, Gerald Tindal, & Brock Rowley", 1) | Out-Null

# ---- Row 5 (SOESD): 11-9-2017 -> 11-8-2018 ----
$t.Cell(5,1).Range.Find.Execute("11-9-2017", $true, $false, $false, $false, $false, $true, 1, $false, "11-8-2018", 1) | Out-Null
$t.Cell(5,2).Range.Paragraphs(1).Range.Find.Execute("Team: Brad Lenhardt & Dan Farley", $true, $false, $false, $false, $false, $true, 1, $false, "Team: Brad Lenhardt & Brock Rowley", 1) | Out-Null

# ---- Row 6 (NWESD): 11-14-2017 -> 11-13-2018 ----
$t.Cell(6,1).Range.Find.Execute("11-14-2017", $true, $false, $false, $false, $false, $true, 1, $false, "11-13-2018", 1) | Out-Null
$t.Cell(6,2).Range.Paragraphs(1).Range.Find.Execute("Team: Brad Lenhardt & Dan Farley", $true, $false, $false, $false, $false, $true, 1, $false, "Team: Brad Lenhardt & Gerald Tindal, Brock Rowley, & Sevrina Tindal", 1) | Out-Null

# ---- Row 7 (Willamette ESD): 11-16-2017 -> 11-15-2018 ----
$t.Cell(7,1).Range.Find.Execute("11-16-2017", $true, $false, $false, $false, $false, $true, 1, $false, "11-15-2018", 1) | Out-Null
$t.Cell(7,2).Range.Paragraphs(1).Range.Find.Execute("Team; Brad Lenhardt & Dan Farley", $true, $false, $false, $false, $false, $true, 1, $false, "Team; Brad Lenhardt, Gerald Tindal, & Brock Rowley", 1) | Out-Null

Write-Output "edits applied"
